# Generate Report for Handback
# Populates the "handback" columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) for the zh-cn and de-de localization status sheets,
# flips the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens a few columns so the longer
# text fits (mirrors what Excel's own column auto-sizing would do after the
# content grew).

$wb = $excel.ActiveWorkbook

$mdDoc1    = "d53922a3-dc98-4b54-88b3-f42ad752343e.md"
$mdDoc2    = "ffff3b0bb8f4-802c-4429-b167-990b74fc2c55.md"
$mdDoc1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c2a7fc138648995beadb5ec0fd5075714bcf6db/e2e/$mdDoc1"
$mdDoc2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c2a7fc138648995beadb5ec0fd5075714bcf6db/e2e/$mdDoc2"

$statusHandedBack = "Handed back: in sync with en-US"

$statusColWidth = 29.9777047293527
$fileColWidth   = 40

$hyperlinkUnderline = 2
$hyperlinkColor     = 6594541

# ---------------------------------------------------------------------------
# Overview sheet: refresh the Status text shown for both localized docs
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# ---------------------------------------------------------------------------
# zh-cn / de-de sheets: stamp the handback info
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Xlf = "d53922a3-dc98-4b54-88b3-f42ad752343e.cb90cff964f0c2913408fef0e4516aabe65151ac.zh-cn.xlf"; HandbackDateTime = "2016-08-17 06:55:53" },
    @{ Name = "de-de"; Xlf = "d53922a3-dc98-4b54-88b3-f42ad752343e.cb90cff964f0c2913408fef0e4516aabe65151ac.de-de.xlf"; HandbackDateTime = "2016-08-17 06:56:00" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C)
    $ws.Range("C2").Value = $statusHandedBack
    $ws.Range("C3").Value = $statusHandedBack

    # Latest Target File (I), Latest Handback File (J), Latest Handback DateTime (K)
    $ws.Range("I2").Value = $mdDoc1
    $ws.Range("J2").Value = $lang.Xlf
    $ws.Range("K2").Value = $lang.HandbackDateTime

    $ws.Range("I3").Value = $mdDoc1
    $ws.Range("J3").Value = $lang.Xlf
    $ws.Range("K3").Value = $lang.HandbackDateTime

    # Rebuild the hyperlinks in display order A2, I2, A3, I3 so relationship
    # ids line up the way a freshly generated report would emit them.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdDoc1Url, [Type]::Missing, [Type]::Missing, $mdDoc1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdDoc1Url, [Type]::Missing, [Type]::Missing, $mdDoc1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdDoc2Url, [Type]::Missing, [Type]::Missing, $mdDoc2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdDoc1Url, [Type]::Missing, [Type]::Missing, $mdDoc1)

    # Style column I like the existing document hyperlinks (underline + blue)
    $ws.Range("I2").Font.Underline = $hyperlinkUnderline
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("I3").Font.Underline = $hyperlinkUnderline
    $ws.Range("I3").Font.Color = $hyperlinkColor

    $ws.Columns.Item(3).ColumnWidth = $statusColWidth
    $ws.Columns.Item(9).ColumnWidth = $fileColWidth
    $ws.Columns.Item(10).ColumnWidth = $fileColWidth
}

Write-Host "Handback report generated for zh-cn and de-de."
